# FSLogix Storage Calculator - update modeled inputs
# Commit: "Add files via upload"
#
# Changes the two primary "Modify:" input cells on the FSLogix Calculator
# sheet (User Count and IOPS Per User). Every other cell in the workbook
# that changes in the diff is a formula that depends on these two values,
# so simply updating the inputs and letting the workbook recalculate
# reproduces the rest of the delta automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FSLogix Calculator")

# Modify: User Count            137 -> 500
$ws.Range("B3").Value = 500

# Modify: IOPS Per User          15 -> 5
$ws.Range("B4").Value = 5

# Update the IOPS data-validation input-message prompt to reflect the new
# suggested values (was "Suggested values:_x000a_Per LoginVSI: 4_x000a_Per
# FSLogix Team: 15_x000a_With Offline Cache and Office365 Cache: 60 - 120").
$ws.Range("B4").Validation.InputMessage = "Suggested values for user profile containers_x000a_Without Office365: 5_x000a_With Office365: 15_x000a_With Offline Cache and Office365 Cache: 60 - 120"

# Restore the active selection on the sheet to D9 (was A14).
$ws.Range("D9").Select() | Out-Null
